$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.846.02"
$ws.Range("E2").Value = "  +1.82%  "

# Row 3
$ws.Range("D3").Value = "1.887.78"
$ws.Range("E3").Value = "  +1.69%  "

# Row 4
$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.74%  "

# Row 5
$ws.Range("D5").Value = "'334.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "

# Row 6
$ws.Range("E6").Value = "  +0.61%  "

# Row 7
$ws.Range("D7").Value = "'0.4732"
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.3931"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.45%  "

# Row 9
$ws.Range("D9").Value = "'47.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.31%  "

# Row 10
$ws.Range("D10").Value = "'0.08066"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.34%  "

# Row 11
$ws.Range("E11").Value = "  +1.29%  "

# Row 12
$ws.Range("D12").Value = "'22.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.84%  "

# Row 13
$ws.Range("D13").Value = "1.892.23"
$ws.Range("E13").Value = "  +2.90%  "

# Row 14
$ws.Range("D14").Value = "'5.992"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.85%  "

# Row 15
$ws.Range("D15").Value = "'7.148"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "

# Row 16
$ws.Range("E16").Value = "  +0.66%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001052"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.02%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06724"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.69%  "

# Row 19
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'87.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.30%  "

# Row 20
$ws.Range("D20").Value = "'17.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.38%  "

# Row 21
$ws.Range("D21").Value = "'1.008"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "

# Row 22
$ws.Range("D22").Value = "27.870.53"
$ws.Range("E22").Value = "  +1.86%  "

# Row 23
$ws.Range("D23").Value = "'5.526"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.74%  "

# Row 24
$ws.Range("D24").Value = "'11.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.71%  "

# Row 25
$ws.Range("D25").Value = "'2.332"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.48%  "

# Row 26
$ws.Range("D26").Value = "2.109.84"
$ws.Range("E26").Value = "  +2.37%  "

# Row 27
$ws.Range("D27").Value = "'159.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.44%  "

# Row 28
$ws.Range("D28").Value = "'20.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.53%  "

# Row 29
$ws.Range("D29").Value = "'2.107"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.88%  "

# Row 30
$ws.Range("D30").Value = "'5.574"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.78%  "

# Row 31
$ws.Range("D31").Value = "'122.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "

# Row 32
$ws.Range("D32").Value = "'0.9786"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.02%  "

# Row 33
$ws.Range("D33").Value = "'0.09497"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "

# Row 34
$ws.Range("D34").Value = "'1.452"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "

# Row 35
$ws.Range("D35").Value = "'3.627"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.06%  "

# Row 36
$ws.Range("D36").Value = "'5.361"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.74%  "

# Row 37
$ws.Range("D37").Value = "'0.06169"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.92%  "

# Row 38
$ws.Range("E38").Value = "  +1.93%  "

# Row 39
$ws.Range("D39").Value = "'1.221"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.56%  "

# Row 40
$ws.Range("D40").Value = "'8.078"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.49%  "

# Row 41
$ws.Range("D41").Value = "'0.6008"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.16%  "

# Row 42
$ws.Range("D42").Value = "'0.1899"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "

# Row 43
$ws.Range("E43").Value = "  +0.89%  "

# Row 44
$ws.Range("E44").Value = "  -1.54%  "

# Row 45
$ws.Range("D45").Value = "'0.5713"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.15%  "

# Row 46
$ws.Range("E46").Value = "  +0.90%  "

# Row 47
$ws.Range("D47").Value = "'3.405"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.57%  "

# Row 48
$ws.Range("E48").Value = "  +0.86%  "

# Row 49
$ws.Range("D49").Value = "'0.06920"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.22%  "

# Row 50
$ws.Range("D50").Value = "'113.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.60%  "

# Row 51
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.073"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.58%  "
